$d = $word.ActiveDocument

# Unicode punctuation helpers
$ldq = [char]0x201c   # left double quotation mark  "
$rsq = [char]0x2019   # right single quotation mark '
$rdq = [char]0x201d   # right double quotation mark "

# ---------------------------------------------------------------
# Change 1: merge the three bold runs of the decision banner into
# a single run: "-- You choose to do enter the fantasy forest --"
# ---------------------------------------------------------------
$search1 = "-- You choose to do enter the fantasy forest --"
$null = $d.Content.Find.Execute($search1, $false, $false, $false, $false, $false, $true, 1, $false, $search1, 2)

# ---------------------------------------------------------------
# Change 2: merge the three Narrator/starting-items runs into one
# run (this also gets rid of the old "_GoBack" bookmark's adjacent
# split, the bookmark itself is removed explicitly afterwards).
# ---------------------------------------------------------------
$search2 = " ${ldq}It${rsq}s been a while since someone chose to go in so quickly, here${rsq}s your starting items,${rdq} Narrator waves her hands and a short sword, a small bag, and a couple of potions materialize onto your hands. You peer into the bag to find a small amount of sparkly gold. `"That bag contains a total of 100 G,${rdq} she states."
$null = $d.Content.Find.Execute($search2, $false, $false, $false, $false, $false, $true, 1, $false, $search2, 2)

# Remove the old "_GoBack" bookmark (it is relocated by change 3 below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------
# Change 3: split "...a strong vicegrip grabs..." into
# "...a strong vice" / "-" / [bookmark _GoBack] / "grip grabs..."
# ---------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("vicegrip")
$splitPos = $idx + 4   # right after "...vice"
$r = $d.Range($splitPos, $splitPos)
$r.InsertAfter("-")

$dashEnd = $splitPos + 1

$tempBefore = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TEMP_BEFORE_DASH", $tempBefore)

$goBackRange = $d.Range($dashEnd, $dashEnd)
$d.Bookmarks.Add("_GoBack", $goBackRange)

$d.Bookmarks("TEMP_BEFORE_DASH").Delete()

# ---------------------------------------------------------------
# Change 4: split the "<MC NAME> ... his story." sentence into
# "You" / " has been gobbled up by the dragon, and the end of " /
# "your" / " story."
# ---------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("<MC NAME>")
$endIdx = $idx + "<MC NAME>".Length
$r = $d.Range($idx, $endIdx)
$r.Text = "You"
$split1 = $idx + 3   # right after "You"

$full = $d.Content.Text
$hisIdx = $full.IndexOf(" his story")
$hisStart = $hisIdx + 1
$hisEnd = $hisStart + 3
$r2 = $d.Range($hisStart, $hisEnd)
$r2.Text = "your"
$split2 = $hisStart          # right before "your"
$split3 = $hisStart + 4      # right after "your"

$bmA = $d.Range($split1, $split1)
$d.Bookmarks.Add("TEMP_SPLIT_A", $bmA)
$bmB = $d.Range($split2, $split2)
$d.Bookmarks.Add("TEMP_SPLIT_B", $bmB)
$bmC = $d.Range($split3, $split3)
$d.Bookmarks.Add("TEMP_SPLIT_C", $bmC)

$d.Bookmarks("TEMP_SPLIT_A").Delete()
$d.Bookmarks("TEMP_SPLIT_B").Delete()
$d.Bookmarks("TEMP_SPLIT_C").Delete()

# ---------------------------------------------------------------
# Change 5: add a new centered "BAD END" paragraph after the
# "Well that was a quick story..." paragraph.
# ---------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Alignment = 1
$newPara.Range.InsertAfter("BAD END")

Write-Output "done"
